$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3500
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -4350
$ws.Range("H32").Value = 16027.526
$ws.Range("I32").Value = 12896.7
$ws.Range("K32").Value = 12896.7
$ws.Range("M32").Value = -12570.7
$ws.Range("H33").Value = 685.6429000000001
$ws.Range("I33").Value = 668.36365
$ws.Range("K33").Value = 668.36365
$ws.Range("M33").Value = -439.36365
$ws.Range("H70").Value = 1540.4546
$ws.Range("I70").Value = 1325.1666
$ws.Range("J70").Value = 1798.8
$ws.Range("K70").Value = 3975.4998
$ws.Range("L70").Value = 5396.4
$ws.Range("M70").Value = -3705.4998
$ws.Range("N70").Value = -5936.4
$ws.Range("H73").Value = 1540.4546
$ws.Range("I73").Value = 1325.1666
$ws.Range("J73").Value = 1798.8
$ws.Range("K73").Value = 3975.4998
$ws.Range("L73").Value = 5396.4
$ws.Range("M73").Value = -3039.4998
$ws.Range("N73").Value = -7268.4
$ws.Range("H74").Value = 12721.833
$ws.Range("J74").Value = 14692.714
$ws.Range("L74").Value = 14692.714
$ws.Range("N74").Value = -16564.714
$ws.Range("H77").Value = 12721.833
$ws.Range("J77").Value = 14692.714
$ws.Range("L77").Value = 73463.57000000001
$ws.Range("N77").Value = -82823.57000000001
$ws.Range("H103").Value = 1824.1428
$ws.Range("I103").Value = 848.75
$ws.Range("J103").Value = 3124.6667
$ws.Range("K103").Value = 2546.25
$ws.Range("L103").Value = 9374.000100000001
$ws.Range("M103").Value = -1960.25
$ws.Range("N103").Value = -10546.0001
$ws.Range("H111").Value = 1253
$ws.Range("I111").Value = 1103.6
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 3310.8
$ws.Range("L111").Value = 6000
$ws.Range("M111").Value = -243.7999999999997
$ws.Range("N111").Value = -12134
$ws.Range("H129").Value = 30316828
$ws.Range("I129").Value = 4077.8
$ws.Range("J129").Value = 55577456
$ws.Range("K129").Value = 12233.4
$ws.Range("L129").Value = 166732368
$ws.Range("M129").Value = -7233.400000000001
$ws.Range("N129").Value = -166742368
$ws.Range("H132").Value = 8828.117
$ws.Range("I132").Value = 6738.533
$ws.Range("J132").Value = 24500
$ws.Range("K132").Value = 20215.599
$ws.Range("L132").Value = 73500
$ws.Range("M132").Value = -17685.599
$ws.Range("N132").Value = -78560
$ws.Range("H137").Value = 6900.161
$ws.Range("I137").Value = 4584.85
$ws.Range("J137").Value = 11109.818
$ws.Range("K137").Value = 13754.55
$ws.Range("L137").Value = 33329.454
$ws.Range("M137").Value = -11204.55
$ws.Range("N137").Value = -38429.454
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3565.2954
$ws.Range("I2").Value = 3091.5134
$ws.Range("J2").Value = 6069.5713
$ws.Range("K2").Value = 3091.5134
$ws.Range("L2").Value = 6069.5713
$ws.Range("M2").Value = -2978.5134
$ws.Range("N2").Value = -6295.5713
$ws.Range("H4").Value = 796.125
$ws.Range("I4").Value = 809.6
$ws.Range("K4").Value = 809.6
$ws.Range("M4").Value = -693.6
$ws.Range("H32").Value = 2175.9656
$ws.Range("I32").Value = 779.1389
$ws.Range("J32").Value = 4461.6816
$ws.Range("K32").Value = 779.1389
$ws.Range("L32").Value = 4461.6816
$ws.Range("M32").Value = -492.1389
$ws.Range("N32").Value = -5035.6816
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36608
$ws.Range("H88").Value = 4637.7646
$ws.Range("I88").Value = 1138.6666
$ws.Range("J88").Value = 5387.5713
$ws.Range("K88").Value = 1138.6666
$ws.Range("L88").Value = 5387.5713
$ws.Range("M88").Value = -732.6666
$ws.Range("N88").Value = -6199.5713
$ws.Range("H91").Value = 4637.7646
$ws.Range("I91").Value = 1138.6666
$ws.Range("J91").Value = 5387.5713
$ws.Range("K91").Value = 1138.6666
$ws.Range("L91").Value = 5387.5713
$ws.Range("M91").Value = 265.3334
$ws.Range("N91").Value = -8195.5713
$ws.Range("H97").Value = 837.7105
$ws.Range("J97").Value = 852.1667
$ws.Range("L97").Value = 852.1667
$ws.Range("N97").Value = -1844.1667
$ws.Range("H116").Value = 3565.2954
$ws.Range("I116").Value = 3091.5134
$ws.Range("J116").Value = 6069.5713
$ws.Range("K116").Value = 3091.5134
$ws.Range("L116").Value = 6069.5713
$ws.Range("M116").Value = -797.5133999999998
$ws.Range("N116").Value = -10657.5713
$ws.Range("H132").Value = 22325.912
$ws.Range("I132").Value = 1579.4
$ws.Range("J132").Value = 61225.625
$ws.Range("K132").Value = 4738.200000000001
$ws.Range("L132").Value = 183676.875
$ws.Range("M132").Value = -2208.200000000001
$ws.Range("N132").Value = -188736.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3565.2954
$ws.Range("I3").Value = 3091.5134
$ws.Range("J3").Value = 6069.5713
$ws.Range("K3").Value = 3091.5134
$ws.Range("L3").Value = 6069.5713
$ws.Range("M3").Value = -2977.5134
$ws.Range("N3").Value = -6297.5713
$ws.Range("H107").Value = 7165.8335
$ws.Range("I107").Value = 6331.6665
$ws.Range("K107").Value = 6331.6665
$ws.Range("M107").Value = -4411.6665
$ws.Range("H133").Value = 100383.75
$ws.Range("J133").Value = 100383.75
$ws.Range("L133").Value = 100383.75
$ws.Range("N133").Value = -110503.75
$ws.Range("H134").Value = 6554.207
$ws.Range("I134").Value = 5844.4165
$ws.Range("J134").Value = 9961.200000000001
$ws.Range("K134").Value = 17533.2495
$ws.Range("L134").Value = 29883.6
$ws.Range("M134").Value = -14998.2495
$ws.Range("N134").Value = -34953.60000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6397.533
$ws.Range("I31").Value = 1910.7142
$ws.Range("J31").Value = 10323.5
$ws.Range("K31").Value = 1910.7142
$ws.Range("L31").Value = 10323.5
$ws.Range("M31").Value = -1615.7142
$ws.Range("N31").Value = -10913.5
$ws.Range("H34").Value = 6397.533
$ws.Range("I34").Value = 1910.7142
$ws.Range("J34").Value = 10323.5
$ws.Range("K34").Value = 1910.7142
$ws.Range("L34").Value = 10323.5
$ws.Range("M34").Value = -1708.7142
$ws.Range("N34").Value = -10727.5
$ws.Range("H99").Value = 4562
$ws.Range("J99").Value = 10000
$ws.Range("L99").Value = 10000
$ws.Range("N99").Value = -12996
$ws.Range("H122").Value = 2421.3333
$ws.Range("I122").Value = 1656
$ws.Range("K122").Value = 4968
$ws.Range("M122").Value = -2518
$ws.Range("H126").Value = 4562
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 5165.6665
$ws.Range("I132").Value = 4460.7617
$ws.Range("K132").Value = 13382.2851
$ws.Range("M132").Value = -10852.2851
$ws.Range("H134").Value = 11465.871
$ws.Range("I134").Value = 10267.223
$ws.Range("J134").Value = 11956.228
$ws.Range("K134").Value = 30801.669
$ws.Range("L134").Value = 35868.68399999999
$ws.Range("M134").Value = -28266.669
$ws.Range("N134").Value = -40938.68399999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6224.6
$ws.Range("I102").Value = 5317.148
$ws.Range("K102").Value = 5317.148
$ws.Range("M102").Value = -3695.148
$ws.Range("H113").Value = 5236.6665
$ws.Range("I113").Value = 1982.2778
$ws.Range("J113").Value = 14999.833
$ws.Range("K113").Value = 1982.2778
$ws.Range("L113").Value = 14999.833
$ws.Range("M113").Value = 187.7221999999999
$ws.Range("N113").Value = -19339.833
$ws.Range("H132").Value = 2486.5454
$ws.Range("I132").Value = 2220.5
$ws.Range("J132").Value = 3196
$ws.Range("K132").Value = 6661.5
$ws.Range("L132").Value = 9588
$ws.Range("M132").Value = -4131.5
$ws.Range("N132").Value = -14648
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H46").Value = 3103.75
$ws.Range("J46").Value = 3425.7144
$ws.Range("L46").Value = 3425.7144
$ws.Range("N46").Value = -3801.7144
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H55").Value = 2192.5151
$ws.Range("I55").Value = 1262.6923
$ws.Range("K55").Value = 1262.6923
$ws.Range("M55").Value = -1089.6923
$ws.Range("H93").Value = 1636.6364
$ws.Range("I93").Value = 3357
$ws.Range("J93").Value = 653.5714
$ws.Range("K93").Value = 3357
$ws.Range("L93").Value = 653.5714
$ws.Range("M93").Value = -2109
$ws.Range("N93").Value = -3149.5714
$ws.Range("H122").Value = 4824.2
$ws.Range("I122").Value = 3335.4
$ws.Range("K122").Value = 10006.2
$ws.Range("M122").Value = -7556.200000000001
$ws.Range("H132").Value = 2268.6086
$ws.Range("I132").Value = 2269.6875
$ws.Range("J132").Value = 2266.1428
$ws.Range("K132").Value = 6809.0625
$ws.Range("L132").Value = 6798.428400000001
$ws.Range("M132").Value = -4279.0625
$ws.Range("N132").Value = -11858.4284
$ws.Range("H136").Value = 4485.8047
$ws.Range("I136").Value = 5184
$ws.Range("J136").Value = 4161.643
$ws.Range("K136").Value = 15552
$ws.Range("L136").Value = 12484.929
$ws.Range("M136").Value = -13002
$ws.Range("N136").Value = -17584.929
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2475.074
$ws.Range("I122").Value = 1919.1818
$ws.Range("K122").Value = 5757.5454
$ws.Range("M122").Value = -3307.5454
$ws.Range("H132").Value = 3672.1968
$ws.Range("I132").Value = 2126.5625
$ws.Range("J132").Value = 9379.154
$ws.Range("K132").Value = 6379.6875
$ws.Range("L132").Value = 28137.462
$ws.Range("M132").Value = -3849.6875
$ws.Range("N132").Value = -33197.462
